# Band2MinPowerInput1_FrequencySweepRanges.xlsx edit
# Change the "float(Power Level)" column (F) on the Band2/Band4/Band5/Band12/Band13
# sheets from 23 to -20, switching the number format to plain "0.0" (no red for
# negatives), extend Band13's table down to row 28 (matching the other band sheets),
# and update the sheet selections / active sheet to reflect the values that were
# last touched while editing.

$wb = $excel.ActiveWorkbook

# ---- Band2, Band4, Band5, Band12: F2:F28, 23 -> -20 --------------------------
$simpleSheets = @("Band2", "Band4", "Band5", "Band12")
foreach ($name in $simpleSheets) {
    $ws = $wb.Worksheets.Item($name)
    $rng = $ws.Range("F2:F28")
    $rng.NumberFormat = "0.0"
    $rng.Value = -20
    $rng.Select()
}

# ---- Band13: F2:F13, 23 -> -20, then extend the table to F28 -----------------
$ws13 = $wb.Worksheets.Item("Band13")
$rng13 = $ws13.Range("F2:F13")
$rng13.NumberFormat = "0.0"
$rng13.Value = -20

# Copy the formatting of the last populated power-level cell down through row 28
# so the newly added rows pick up the same style (border/fill/font/number format).
$ws13.Range("F13").Copy()
$ws13.Range("F14:F28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws13.Range("F14:F28").Select()

# ---- Band66: becomes the active sheet/tab, with G23 selected -----------------
$ws66 = $wb.Worksheets.Item("Band66")
$ws66.Range("G23").Select()
